# Updated legacy GSC export data:
# The daily export table on the "Chart" sheet is a rolling window of dates
# (one row per day, oldest first). The refresh drops the oldest date
# (2025-10-16) and appends the next day (2026-01-14) with zero counts,
# which shifts every other row up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest day's row -> every subsequent row (dates + counts)
# shifts up by one, matching the diff's row-2..row-11 value changes.
$ws.Rows.Item(2).Delete()

# Append the new day as the new last row (row 91). Stage the date text in
# a scratch cell formatted as Text first, then copy just the *value* into
# place - otherwise Excel's normal text entry would parse a "YYYY-MM-DD"
# looking string as a date serial instead of keeping it as literal text
# (which is how every other date in this column is stored).
$scratch = $ws.Range("E1")
$scratch.NumberFormat = "@"
$scratch.Value = "2026-01-14"
$scratch.Copy()
$ws.Range("A91").PasteSpecial(-4163)   # xlPasteValues
$scratch.EntireColumn.Delete()         # remove the scratch column entirely

# New day has no recorded URLs yet, consistent with the rest of the tail.
$ws.Range("B91").Value = 0
$ws.Range("C91").Value = 0

Write-Output ("Chart used range: " + $ws.UsedRange.Address())
Write-Output ("A2=" + $ws.Range("A2").Value() + " C2=" + $ws.Range("C2").Value())
Write-Output ("A91=" + $ws.Range("A91").Value() + " B91=" + $ws.Range("B91").Value() + " C91=" + $ws.Range("C91").Value())
